# Update "想去人数" (number of people interested) values in the
# "展览" and "全部类型" sheets, rows 3-6, column F.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 97
    $ws.Range("F4").Value = 80
    $ws.Range("F5").Value = 2537
    $ws.Range("F6").Value = 233
}
